# Applies the "EC" (Estado de Cuenta) update:
#  - Shrinks the worker table from 11 rows down to 2 rows (deleting the
#    JAVIER/JOSHUAN rows and the duplicated period rows), keeping only
#    MARCO FIDEL GOMEZ ARRIETA and ABEL ANTONIO OROZCO TEHERAN.
#  - Updates the Valor Mora / Cant. Trabajadores / Cant. Periodos summary
#    figures to match the new (smaller) data set.
#  - Updates the Valor Mora / Salario Basico figures on the two remaining
#    worker rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Summary figures (rows 11 & 13) -------------------------------------
$ws.Range("E11").Value = 35200
$ws.Range("C13").Value = 2
$ws.Range("F13").Value = 2

# --- Rewrite the two worker rows that survive (currently rows 16 & 17) --
$ws.Range("B16").Value = "CC"
$ws.Range("C16").Value = "6866376"
$ws.Range("D16").Value = "MARCO FIDEL GOMEZ ARRIETA"
$ws.Range("E16").Value = "1705"
$ws.Range("F16").Value = 25600
$ws.Range("G16").Value = 1300000

$ws.Range("B17").Value = "CC"
$ws.Range("C17").Value = "1137220356"
$ws.Range("D17").Value = "ABEL ANTONIO OROZCO TEHERAN"
$ws.Range("E17").Value = "1907"
$ws.Range("F17").Value = 9600
$ws.Range("G17").Value = 1200000

# --- Remove the now-obsolete rows (old rows 18-26) ----------------------
# This shifts the trailing signature-block rows (old 31/32) up to 22/23.
$ws.Rows("18:26").Delete()
